# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker ABELARDO ENRIQUE VALDES DIAZ (1052096935) is removed entirely,
# and new periods (2307-2312, 2401-2412) are added for the two remaining
# workers (ABEL ANTONIO HIGUITA GUISAO / 71252818 and
# JOSE JULIAN GIRALDO GOEZ / 1001669965). The summary totals and the table
# footer shift up accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Shrink the data block by 4 rows -----------------------------------
# Deleting from the TOP of the data block (rather than the tail) lets the
# specially-formatted "closing" row (originally row 42, with its heavier
# bottom border) slide up and land on the new last data row (38) with its
# formatting intact, while every other data row keeps the regular style.
$ws.Rows("17:20").Delete()

# --- 2) Rewrite the data table (rows 16-38) --------------------------------
$data = New-Object 'object[,]' 23,6
$rows = @(
  @("CC","71252818","ABEL ANTONIO HIGUITA GUISAO","2307",56000,3000000),
  @("CC","71252818","ABEL ANTONIO HIGUITA GUISAO","2308",120000,3000000),
  @("CC","71252818","ABEL ANTONIO HIGUITA GUISAO","2309",120000,3000000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2309",76266,2200000),
  @("CC","71252818","ABEL ANTONIO HIGUITA GUISAO","2310",120000,3000000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2310",88000,2200000),
  @("CC","71252818","ABEL ANTONIO HIGUITA GUISAO","2311",120000,3000000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2311",88000,2200000),
  @("CC","71252818","ABEL ANTONIO HIGUITA GUISAO","2312",120000,3000000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2312",88000,2200000),
  @("CC","71252818","ABEL ANTONIO HIGUITA GUISAO","2401",120000,3000000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2401",88000,2200000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2402",88000,2200000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2403",88000,2200000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2404",88000,2200000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2405",88000,2200000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2406",88000,2200000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2407",88000,2200000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2408",88000,2200000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2409",88000,2200000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2410",88000,2200000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2411",88000,2200000),
  @("CC","1001669965","JOSE JULIAN GIRALDO GOEZ","2412",88000,2200000)
)
for ($i = 0; $i -lt $rows.Count; $i++) {
  for ($j = 0; $j -lt 6; $j++) {
    $data[$i,$j] = $rows[$i][$j]
  }
}
$ws.Range("B16:G38").Value = $data

# --- 3) Update the summary figures -----------------------------------------
$ws.Range("E11").Value = 2172266          # VALOR MORA (sum of Valor Mora column)
$ws.Range("C13").Value = 2                # Cant. Trabajadores
$ws.Range("F13").Value = 18               # Cant. Periodos

# --- 4) Column D ("Nombre Trabajador") is now narrower now that the longest
# name in the table is shorter than before -----------------------------
$ws.Columns("D").ColumnWidth = 29.54296875
